$d = $word.ActiveDocument

# Robustly locate the last content paragraph (the one that ends with the
# BeautifulSoup "print(tag.get('href', None))" example followed by an
# ellipsis) by scanning the Paragraphs collection for unique text, rather
# than relying on a hard-coded paragraph index.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*print(tag.get*") {
        $targetIndex = $i
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate anchor paragraph containing the BeautifulSoup example."
}

$anchorPara = $d.Paragraphs.Item($targetIndex)
$rng = $anchorPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1).Range

$bodyFragment = '<w:p><w:r><w:t xml:space="preserve">XML - </w:t></w:r><w:r><w:t>define and store data in a shareable manner</w:t></w:r><w:r><w:t>. Has a Start tag,end tag, text content, attribute, self closing tag.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We have a built in XML parser in python called </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>element tree</w:t></w:r><w:r><w:t xml:space="preserve">. Ie.. </w:t></w:r><w:r><w:t>import xml.etree.ElementTree as ET</w:t></w:r><w:r><w:t>.   (here as ET is just a shortcut handle to call the element tree)</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Tree =  ET.fromstring(data) this will read this data from string and give us an object(tree).</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We can find using tags so using uaing object.find we can find the tag name. ie.. </w:t></w:r><w:r><w:t>print(''Name:'', tree.find(''name'').text)</w:t></w:r><w:r><w:t>. to get the contents of an attribute we use .get ie..</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>print(''Attr:'', tree.find(''email'').get(''hide''))</w:t></w:r></w:p><w:p><w:r><w:t>So</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> to get attribute we use get(fn) to get the info like name age that we provide which is not an attribute we use find(fn).</w:t></w:r><w:r><w:t xml:space="preserve"> Findall – to get all the users and find to just get a single user.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">JSON – represent data as nested “lists” and “dictionaries”. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Here we take the string from our given input(ie data) then pass it into json library loads which reads it, parses it and looks at all the white spaces and returns us a dictionary which we assign an object like xml we did before ie info. Eg </w:t></w:r><w:r><w:t>print(''Name:'', info["name"])</w:t></w:r><w:r><w:t>. Gives chuck as output.</w:t></w:r></w:p><w:p><w:r><w:t>API – Application program interface. They set the rules, set the urls. Will say if its xml or json.</w:t></w:r></w:p><w:p><w:r><w:t>Objects – are bits of code and data. Object hides details, they allow the rest of the program to ignore the details about “us”.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Class – a template  - eg. </w:t></w:r><w:r><w:t>D</w:t></w:r><w:r><w:t>og</w:t></w:r><w:r><w:t xml:space="preserve"> or a cookie cuter</w:t></w:r></w:p><w:p><w:r><w:t>Method or message – A defined capability of a class – eg. bark()</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Field or attribute – a bit of data in a class – eg. length </w:t></w:r></w:p><w:p><w:r><w:t>Object or instance – a particular instance of a class  - eg. Bulldog</w:t></w:r><w:r><w:t xml:space="preserve"> or snoman cookie, apple cookie(BASICALLY MEAN THE SHAPE HERE) </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Code we write is in the class and object are the instance. Methods are fns are lives inside the class. Fields are the variables that are defined inside the class. </w:t></w:r></w:p><w:p><w:r><w:t>Dir() – is used to find the capabilities of our newly created class.</w:t></w:r></w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.InsertXML($packageXml)

Write-Host "Paragraphs count: $($d.Paragraphs.Count)"
